$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "names"
$ws.Range("B1").Value = "salarys"

# Row 2 - ali
$ws.Range("A2").Value = "ali "
$ws.Range("B2").Value = 3000

# Row 3 - kemal
$ws.Range("A3").Value = "kemal"
$ws.Range("B3").Value = 40000

# Currency column filled last (matches shared-string build order: tl follows kemal)
$ws.Range("C2").Value = "tl"
$ws.Range("C3").Value = "tl"

# Leave the active cell where data entry ended, just below the last row
$ws.Range("C4").Select() | Out-Null
